# Update QS doc example
# Insert a new "Building" twin row above the existing Floor/Room rows so the
# sample shows a Building -> Floor -> Room relationship hierarchy.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new Building row right under the header row.
$ws.Rows("2:2").Insert()

# New row 2: the Building twin itself.
$ws.Range("A2").Value2 = "dtmi:example:Building;1"
$ws.Range("B2").Value2 = "BuildingA"
$ws.Range("E2").Value2 = '{"Date": 2001}'

# The two Floor rows (now rows 3 & 4) relate to the new Building via a
# "contains" relationship.
$ws.Range("C3").Value2 = "BuildingA"
$ws.Range("D3").Value2 = "contains"

$ws.Range("C4").Value2 = "BuildingA"
$ws.Range("D4").Value2 = "contains"

# Narrow column A now that it no longer needs to fit the longest label.
$ws.Columns("A").ColumnWidth = 26.14

# Update the active selection / scroll position like a user would after
# finishing the edit.
$ws.Range("D4").Select()
